$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.25
$ws.Range("D3").Value = 1.43
$ws.Range("B4").Value = 1.49
$ws.Range("C4").Value = 1.42
$ws.Range("E4").Value = 1.26
$ws.Range("F4").Value = 1.06
$ws.Range("D6").Value = 1.55
$ws.Range("G6").Value = 1.01
$ws.Range("F7").Value = 1.54
